$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update financial figures for the most recent period (column D)
$ws.Range("D45").Value = 145600    # Other Current Assets
$ws.Range("D46").Value = 637900    # Total Current Assets

# Long Term Investments row was "NA" text across D47:J47; now numeric 0
$ws.Range("D47:J47").Value = 0

$ws.Range("D48").Value = 10874800  # Property Plant and Equipment
$ws.Range("D52").Value = 213200    # Other Assets
$ws.Range("D54").Value = 9209900   # Total Assets
$ws.Range("D59").Value = 866200    # Other Current Liabilities
$ws.Range("D60").Value = 616700    # Total Current Liabilities
$ws.Range("D62").Value = 1281200   # Other Liabilities
$ws.Range("D66").Value = 7709600   # Total Liabilities
$ws.Range("D72").Value = -366000   # Retained Earnings
$ws.Range("D76").Value = 1500300   # Total Stockholder Equity
